$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (price + 1h volume change).
# Cells are stored as text, so force Text number format before
# assigning to avoid Excel auto-converting to numeric/percentage types
# (which would drop trailing zeros / reformat the percent sign).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.50%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.12%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.124"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.85%"

# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.39%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.275"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.65%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.617"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.09%"

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.68%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9092"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.83%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1180"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.55%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1802"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.15%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09218"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.54%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04264"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.91%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1043"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.10%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001252"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.27%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005861"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.83%"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.18%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.918"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.36%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1369"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2736"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.16%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04059"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.23%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001271"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.56%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004078"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.74%"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.38%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003745"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02431"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.92%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05262"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.64%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007799"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.85%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1303"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.09%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006782"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.20%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001949"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.06%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007547"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.76%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3079"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.34%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006897"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.21%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08051"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,666.32%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.05%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
